$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the Work Items rows (A:H) -------------------------------------
# Row 2: Requisitos do Sistema de Rastreamento (text unchanged, values unchanged)
$ws.Range("A2").Value = "Requisitos do Sistema de Rastreamento"
$ws.Range("B2").Value = "Alta"
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = 1
$ws.Range("G2").Value = 240
$ws.Range("H2").Value = 1.6

# Row 3: Realizar Análise de Sistema (text unchanged, values unchanged)
$ws.Range("A3").Value = "Realizar Análise de Sistema"
$ws.Range("B3").Value = "Alta"
$ws.Range("C3").Value = 15
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 360
$ws.Range("H3").Value = 2.1

# Row 4: Criar Projeto do Sistema de Rastreamento (text unchanged, values unchanged)
$ws.Range("A4").Value = "Criar Projeto do Sistema de Rastreamento"
$ws.Range("B4").Value = "Alta"
$ws.Range("C4").Value = 40
$ws.Range("D4").Value = 1
$ws.Range("G4").Value = 480
$ws.Range("H4").Value = 3.25

# Row 5: Implementação do Sistema de Rastreamento (text unchanged, values unchanged)
$ws.Range("A5").Value = "Implementação do Sistema de Rastreamento"
$ws.Range("B5").Value = "Alta"
$ws.Range("C5").Value = 20
$ws.Range("D5").Value = 0.35
$ws.Range("G5").Value = 530
$ws.Range("H5").Value = 2.2

# Row 6: now "Criar tela de Login do Sistema de Rastreamento" (was "Realizar Testes no Sistema de Rastreamento")
$ws.Range("A6").Value = "Criar tela de Login do Sistema de Rastreamento"
$ws.Range("B6").Value = "Alta"
$ws.Range("C6").Value = 40
$ws.Range("D6").Value = 1
$ws.Range("G6").Value = 72
$ws.Range("H6").Value = 0.83

# Row 9 text is set before row 7's so the two brand-new shared strings are
# interned in the same order the original authors typed them in (keeps the
# shared-string table layout byte-identical to the source file).
# Row 9: now "Realizar Testes Tela Localização de veículos" (was "Criar Tela de  Localização de Veículos")
$ws.Range("A9").Value = "Realizar Testes Tela Localização de veículos"
$ws.Range("B9").Value = "Média"
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = 0.45
$ws.Range("G9").Value = 160
$ws.Range("H9").Value = 0.52

# Row 7: now "Realizar Testes Tela Login no Sistema de Ratsreamento" (was "Implantação do Sistema de Rastreamento")
$ws.Range("A7").Value = "Realizar Testes Tela Login no Sistema de Ratsreamento"
$ws.Range("B7").Value = "Média"
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = 0.84
$ws.Range("G7").Value = 160
$ws.Range("H7").Value = 1.2

# Row 8: now "Criar Tela de  Localização de Veículos" (was "Criar tela de Login do Sistema de Rastreamento")
$ws.Range("A8").Value = "Criar Tela de  Localização de Veículos"
$ws.Range("B8").Value = "Alta"
$ws.Range("C8").Value = 60
$ws.Range("D8").Value = 0.4
$ws.Range("G8").Value = 120
$ws.Range("H8").Value = 0.23

# --- Remove two now-unused blank rows so the trailing formatted row shifts
#     up from row 25 to row 23 (dimension A1:I25 -> A1:I23) ----------------
$ws.Rows("10:11").Delete()

# --- Column A got a bit wider -------------------------------------------
$ws.Columns("A").ColumnWidth = 53

# --- Selection moved from A14 to A13 -------------------------------------
$ws.Range("A13").Select() | Out-Null
